$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Design Info")

# Remove the blank top row so the content block shifts up by one row.
$ws.Rows("1:1").Delete()

# Remove the blank row that used to sit between the old rows 49 and 51
# (now row 49 after the first deletion), tightening the design-info block.
$ws.Rows("49:49").Delete()

$ws.Activate()
$ws.Range("D4").Select()
